$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on humidity (%) cells so Excel does not auto-convert
# the literal strings like "90%" into a numeric 0.9 percentage value.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H36").NumberFormat = "@"

$ws.Range('E2').Value = '2026-02-06 14:47:41'
$ws.Range('H2').Value = '90%'
$ws.Range('K2').Value = '7.9 MJ/m2'
$ws.Range('M2').Value = '5.8 °C 14:06 TU'
$ws.Range('O2').Value = '-0.4 °C'
$ws.Range('E3').Value = '2026-02-06 14:47:44'
$ws.Range('H3').Value = '70%'
$ws.Range('K3').Value = '10.7 MJ/m2'
$ws.Range('O3').Value = '-2.0 °C'
$ws.Range('E4').Value = '2026-02-06 14:47:46'
$ws.Range('H4').Value = '57%'
$ws.Range('J4').Value = '996.5 hPa'
$ws.Range('K4').Value = '10.7 MJ/m2'
$ws.Range('O4').Value = '13.3 °C'
$ws.Range('E5').Value = '2026-02-06 14:47:48'
$ws.Range('H5').Value = '69%'
$ws.Range('J5').Value = '996.8 hPa'
$ws.Range('K5').Value = '9.6 MJ/m2'
$ws.Range('O5').Value = '10.5 °C'
$ws.Range('E6').Value = '2026-02-06 14:47:51'
$ws.Range('J6').Value = '998.0 hPa'
$ws.Range('K6').Value = '8.6 MJ/m2'
$ws.Range('O6').Value = '15.3 °C'
$ws.Range('E7').Value = '2026-02-06 14:47:53'
$ws.Range('H7').Value = '62%'
$ws.Range('J7').Value = '997.6 hPa'
$ws.Range('K7').Value = '11.2 MJ/m2'
$ws.Range('O7').Value = '11.4 °C'
$ws.Range('E8').Value = '2026-02-06 14:47:55'
$ws.Range('H8').Value = '78%'
$ws.Range('K8').Value = '10.7 MJ/m2'
$ws.Range('O8').Value = '9.6 °C'
$ws.Range('E9').Value = '2026-02-06 14:47:58'
$ws.Range('H9').Value = '87%'
$ws.Range('O9').Value = '4.1 °C'
$ws.Range('E10').Value = '2026-02-06 14:48:00'
$ws.Range('H10').Value = '89%'
$ws.Range('O10').Value = '8.3 °C'
$ws.Range('E11').Value = '2026-02-06 14:48:02'
$ws.Range('H11').Value = '79%'
$ws.Range('K11').Value = '7.6 MJ/m2'
$ws.Range('O11').Value = '4.9 °C'
$ws.Range('E12').Value = '2026-02-06 14:48:05'
$ws.Range('K12').Value = '10.7 MJ/m2'
$ws.Range('O12').Value = '14.2 °C'
$ws.Range('E13').Value = '2026-02-06 14:48:07'
$ws.Range('H13').Value = '78%'
$ws.Range('O13').Value = '9.6 °C'
$ws.Range('E14').Value = '2026-02-06 14:48:09'
$ws.Range('K14').Value = '6.4 MJ/m2'
$ws.Range('E15').Value = '2026-02-06 14:48:11'
$ws.Range('H15').Value = '74%'
$ws.Range('J15').Value = '996.9 hPa'
$ws.Range('K15').Value = '10.6 MJ/m2'
$ws.Range('O15').Value = '9.8 °C'
$ws.Range('E16').Value = '2026-02-06 14:48:14'
$ws.Range('H16').Value = '87%'
$ws.Range('K16').Value = '8.3 MJ/m2'
$ws.Range('O16').Value = '5.5 °C'
$ws.Range('E17').Value = '2026-02-06 14:48:16'
$ws.Range('H17').Value = '88%'
$ws.Range('K17').Value = '9.0 MJ/m2'
$ws.Range('L17').Value = '16.6 km/h - 255º 14:26 TU'
$ws.Range('M17').Value = '12.5 °C 14:28 TU'
$ws.Range('O17').Value = '5.3 °C'
$ws.Range('E18').Value = '2026-02-06 14:48:19'
$ws.Range('K18').Value = '5.1 MJ/m2'
$ws.Range('E19').Value = '2026-02-06 14:48:21'
$ws.Range('H19').Value = '79%'
$ws.Range('K19').Value = '10.3 MJ/m2'
$ws.Range('O19').Value = '9.3 °C'
$ws.Range('E20').Value = '2026-02-06 14:48:24'
$ws.Range('K20').Value = '10.4 MJ/m2'
$ws.Range('O20').Value = '-1.9 °C'
$ws.Range('E21').Value = '2026-02-06 14:48:26'
$ws.Range('K21').Value = '9.5 MJ/m2'
$ws.Range('O21').Value = '7.9 °C'
$ws.Range('E22').Value = '2026-02-06 14:48:28'
$ws.Range('K22').Value = '10.1 MJ/m2'
$ws.Range('O22').Value = '10.2 °C'
$ws.Range('E23').Value = '2026-02-06 14:48:30'
$ws.Range('J23').Value = '996.9 hPa'
$ws.Range('K23').Value = '8.5 MJ/m2'
$ws.Range('L23').Value = '21.6 km/h - 304º 14:05 TU'
$ws.Range('O23').Value = '9.7 °C'
$ws.Range('E24').Value = '2026-02-06 14:48:33'
$ws.Range('J24').Value = '996.3 hPa'
$ws.Range('K24').Value = '10.6 MJ/m2'
$ws.Range('M24').Value = '16.4 °C 14:14 TU'
$ws.Range('O24').Value = '13.0 °C'
$ws.Range('E25').Value = '2026-02-06 14:48:35'
$ws.Range('H25').Value = '82%'
$ws.Range('I25').Value = '0.2 mm'
$ws.Range('J25').Value = '998.1 hPa'
$ws.Range('K25').Value = '8.0 MJ/m2'
$ws.Range('L25').Value = '20.9 km/h - 239º 14:17 TU'
$ws.Range('M25').Value = '10.4 °C 14:27 TU'
$ws.Range('O25').Value = '3.8 °C'
$ws.Range('E26').Value = '2026-02-06 14:48:38'
$ws.Range('K26').Value = '7.4 MJ/m2'
$ws.Range('E27').Value = '2026-02-06 14:48:40'
$ws.Range('H27').Value = '84%'
$ws.Range('J27').Value = '996.9 hPa'
$ws.Range('K27').Value = '9.5 MJ/m2'
$ws.Range('L27').Value = '21.2 km/h - 136º 14:09 TU'
$ws.Range('O27').Value = '10.2 °C'
$ws.Range('E28').Value = '2026-02-06 14:48:42'
$ws.Range('H28').Value = '85%'
$ws.Range('O28').Value = '4.2 °C'
$ws.Range('E29').Value = '2026-02-06 14:48:45'
$ws.Range('H29').Value = '61%'
$ws.Range('K29').Value = '11.2 MJ/m2'
$ws.Range('O29').Value = '12.3 °C'
$ws.Range('E30').Value = '2026-02-06 14:48:47'
$ws.Range('K30').Value = '8.3 MJ/m2'
$ws.Range('E31').Value = '2026-02-06 14:48:50'
$ws.Range('H31').Value = '88%'
$ws.Range('O31').Value = '6.7 °C'
$ws.Range('E32').Value = '2026-02-06 14:48:52'
$ws.Range('J32').Value = '998.4 hPa'
$ws.Range('K32').Value = '10.7 MJ/m2'
$ws.Range('O32').Value = '15.9 °C'
$ws.Range('E33').Value = '2026-02-06 14:48:55'
$ws.Range('O33').Value = '9.6 °C'
$ws.Range('E34').Value = '2026-02-06 14:48:57'
$ws.Range('H34').Value = '77%'
$ws.Range('K34').Value = '10.5 MJ/m2'
$ws.Range('O34').Value = '8.2 °C'
$ws.Range('E35').Value = '2026-02-06 14:49:00'
$ws.Range('K35').Value = '7.9 MJ/m2'
$ws.Range('O35').Value = '-2.2 °C'
$ws.Range('E36').Value = '2026-02-06 14:49:02'
$ws.Range('H36').Value = '61%'
$ws.Range('K36').Value = '10.5 MJ/m2'
